# peakList accessors for MSPeakListsSet
# Mark the newly-added "done" (column G) / "ionize" (column F) support cells
# with "X" on the "mslists" sheet for the peakList-related accessor rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("mslists")

# row 7  - averagedPeakLists: add ionize (F) + done (G)
$ws.Range("F7").Value = "X"
$ws.Range("G7").Value = "X"

# row 10 - groupNames: add done (G)
$ws.Range("G10").Value = "X"

# row 11 - groupFeatIndex: add done (G)
$ws.Range("G11").Value = "X"

# row 12 - groupAlgorithm: add done (G)
$ws.Range("G12").Value = "X"

# row 13 - peakLists: add ionize (F) + done (G)
$ws.Range("F13").Value = "X"
$ws.Range("G13").Value = "X"

# Reflect the last-selected cell in the saved view state
$ws.Range("G15").Select()
